# Apply the "soldes" (leave balance) feature: add a Site/GRH lookup column to
# the Managers sheet, and add Site / Congé N / Congé N-1 columns to the
# Salariés sheet, then leave the UI focused on the Salariés sheet.

$wb = $excel.ActiveWorkbook

$wsGRH      = $wb.Worksheets.Item("GRH")
$wsManagers = $wb.Worksheets.Item("Managers")
$wsSalaries = $wb.Worksheets.Item("Salariés")

# ---------------------------------------------------------------------------
# Managers sheet: new column C = "GRH" code the manager belongs to.
# ---------------------------------------------------------------------------

# Borrow the bordered style used on the GRH sheet's header cell so the new
# column matches the sheet's existing look (style index 2) rather than the
# Managers sheet's own A/B style (index 3).
$wsGRH.Range("A1").Copy()
$wsManagers.Range("C1:C3").PasteSpecial(-4122)

$wsManagers.Range("C1").Value = "GRH"
$wsManagers.Range("C2").Value = "G001"
$wsManagers.Range("C3").Value = "G002"

# ---------------------------------------------------------------------------
# Salariés sheet: new columns F (Congé N), G (Congé N-1), E (Site).
# Column order of authoring matches the order new shared strings appear in.
# ---------------------------------------------------------------------------

# Apply the existing header/body style (index 2, same as columns A-D) to the
# new F/G/E columns in one shot via copy/paste of formats.
$wsSalaries.Range("A1:A21").Copy()
$wsSalaries.Range("E1:G21").PasteSpecial(-4122)

$wsSalaries.Range("F1").Value = "Congé N"
$wsSalaries.Range("F2").Value = 5
$wsSalaries.Range("F3").Value = 5
$wsSalaries.Range("F4").Value = 5
$wsSalaries.Range("F5").Value = 5
$wsSalaries.Range("F6").Value = 0
$wsSalaries.Range("F7").Value = 0
$wsSalaries.Range("F8").Value = 0
$wsSalaries.Range("F9").Value = 2.5
$wsSalaries.Range("F10").Value = 2.5
$wsSalaries.Range("F11").Value = 2.5
$wsSalaries.Range("F12").Value = 2.5
$wsSalaries.Range("F13").Value = 2.5
$wsSalaries.Range("F14").Value = 5
$wsSalaries.Range("F15").Value = 5
$wsSalaries.Range("F16").Value = 5
$wsSalaries.Range("F17").Value = 5
$wsSalaries.Range("F18").Value = 5
$wsSalaries.Range("F19").Value = 5
$wsSalaries.Range("F20").Value = 5
$wsSalaries.Range("F21").Value = 5

$wsSalaries.Range("G1").Value = "Congé N-1"
$wsSalaries.Range("G2").Value = 22
$wsSalaries.Range("G3").Value = 14
$wsSalaries.Range("G4").Value = 18
$wsSalaries.Range("G5").Value = 20
$wsSalaries.Range("G6").Value = 11
$wsSalaries.Range("G7").Value = 30
$wsSalaries.Range("G8").Value = 30
$wsSalaries.Range("G9").Value = 30
$wsSalaries.Range("G10").Value = 25
$wsSalaries.Range("G11").Value = 20
$wsSalaries.Range("G12").Value = 10
$wsSalaries.Range("G13").Value = 7
$wsSalaries.Range("G14").Value = 17
$wsSalaries.Range("G15").Value = 20
$wsSalaries.Range("G16").Value = 20
$wsSalaries.Range("G17").Value = 30
$wsSalaries.Range("G18").Value = 30
$wsSalaries.Range("G19").Value = 30
$wsSalaries.Range("G20").Value = 14
$wsSalaries.Range("G21").Value = 19

$wsSalaries.Range("E1").Value = "Site"
$wsSalaries.Range("E2").Value = "Site 01"
$wsSalaries.Range("E3").Value = "Site 02"
$wsSalaries.Range("E4").Value = "Site 03"
$wsSalaries.Range("E5").Value = "Site 04"
$wsSalaries.Range("E6").Value = "Site 05"
$wsSalaries.Range("E7").Value = "Site 01"
$wsSalaries.Range("E8").Value = "Site 02"
$wsSalaries.Range("E9").Value = "Site 03"
$wsSalaries.Range("E10").Value = "Site 04"
$wsSalaries.Range("E11").Value = "Site 05"
$wsSalaries.Range("E12").Value = "Site 01"
$wsSalaries.Range("E13").Value = "Site 02"
$wsSalaries.Range("E14").Value = "Site 03"
$wsSalaries.Range("E15").Value = "Site 04"
$wsSalaries.Range("E16").Value = "Site 05"
$wsSalaries.Range("E17").Value = "Site 01"
$wsSalaries.Range("E18").Value = "Site 02"
$wsSalaries.Range("E19").Value = "Site 03"
$wsSalaries.Range("E20").Value = "Site 04"
$wsSalaries.Range("E21").Value = "Site 05"

# New F/G columns get a custom width matching column B's.
$wsSalaries.Range("F1:G21").ColumnWidth = 12.93

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping (mirrors what Excel records when a
# user clicks around while editing): Managers sheet ends with C17:C18
# selected, and the Salariés sheet is the final active tab with I9 selected.
# ---------------------------------------------------------------------------

$wsManagers.Range("C17:C18").Select()
$wsSalaries.Range("I9").Select()
